$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark supply items (rows 10-14) as completed: bump the "E" quantity column
# to 4 and flip the "G" completion flag from 0 to 1.
10..14 | ForEach-Object {
    $row = $_
    $ws.Cells.Item($row, 5).Value = 4   # column E
    $ws.Cells.Item($row, 7).Value = 1   # column G
}

# Reflect the in-progress work location: scroll the sheet down slightly and
# leave the selection on G10.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G10").Select()
